$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column C: header + values (mapping of HEBT CEB slots)
$ws.Range("C1").Value = "HEBT_CEB"
$ws.Range("C2").Value = "SPARE"
$ws.Range("C3").Value = "H5_001B_CEB"
$ws.Range("C4").Value = "U1_023A_CEB"
$ws.Range("C5").Value = "H2_019A_CEB"
$ws.Range("C7").Value = "V1_044A_CEB"
$ws.Range("C8").Value = "T2_015A_CEB"
$ws.Range("C9").Value = "V1_030A_CEB"
$ws.Range("C10").Value = "HE_027A_CEB"
$ws.Range("C11").Value = "V1_003A_CEB"
$ws.Range("C12").Value = "T1_011A_CEB"
$ws.Range("C13").Value = "H4_016A_CEB"
$ws.Range("C14").Value = "Z1_011A_CEB"
$ws.Range("C15").Value = "U1_003A_CEB"
$ws.Range("C16").Value = "Z2_008A_CEB"
$ws.Range("C17").Value = "H5_012A_CEB"
$ws.Range("C18").Value = "T2_008A_CEB"
$ws.Range("C19").Value = "V2_013A_CEB"
$ws.Range("C20").Value = "Z2_015A_CEB"
$ws.Range("C21").Value = "U1_011A_CEB"
$ws.Range("C22").Value = "U2_013A_CEB"

# Header formatting: bold, matching style of A1/B1
$ws.Range("C1").Font.Bold = $true

# Column width best fit
$ws.Columns("C").EntireColumn.AutoFit()

# Update the absolute path recorded by Excel (x15ac:absPath)
$wb.Path = "D:\optics_updateMapping\"

# Update the active selection to C1
$ws.Range("C1").Select()
